$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 18066318
$ws.Range("I33").Value = 10462092
$ws.Range("K33").Value = 10462092
$ws.Range("M33").Value = -10461863
$ws.Range("H53").Value = 41667170
$ws.Range("I53").Value = 519
$ws.Range("K53").Value = 519
$ws.Range("M53").Value = 118
$ws.Range("H70").Value = 113232.89
$ws.Range("J70").Value = 2616.3333
$ws.Range("L70").Value = 7848.999899999999
$ws.Range("N70").Value = -8388.999899999999
$ws.Range("H73").Value = 113232.89
$ws.Range("J73").Value = 2616.3333
$ws.Range("L73").Value = 7848.999899999999
$ws.Range("N73").Value = -9720.999899999999
$ws.Range("H74").Value = 8754.241
$ws.Range("I74").Value = 8754.241
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 8754.241
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -7818.241
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 8754.241
$ws.Range("I77").Value = 8754.241
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 43771.205
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -39091.205
$ws.Range("N77").Value = ""
$ws.Range("H113").Value = 5682.7744
$ws.Range("I113").Value = 4086.1177
$ws.Range("J113").Value = 7621.5713
$ws.Range("K113").Value = 4086.1177
$ws.Range("L113").Value = 7621.5713
$ws.Range("M113").Value = -832.1176999999998
$ws.Range("N113").Value = -14129.5713
$ws.Range("H118").Value = 591.3
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
$ws.Range("H125").Value = 456152.28
$ws.Range("I125").Value = 1421516.5
$ws.Range("J125").Value = 1863.2354
$ws.Range("K125").Value = 12793648.5
$ws.Range("L125").Value = 16769.1186
$ws.Range("M125").Value = -12791188.5
$ws.Range("N125").Value = -21689.1186
$ws.Range("H127").Value = 2974.3333
$ws.Range("I127").Value = 1569.2
$ws.Range("J127").Value = 10000
$ws.Range("K127").Value = 4707.6
$ws.Range("L127").Value = 30000
$ws.Range("M127").Value = 252.3999999999996
$ws.Range("N127").Value = -39920
$ws.Range("H137").Value = 1483.0377
$ws.Range("I137").Value = 884.65216
$ws.Range("K137").Value = 2653.95648
$ws.Range("M137").Value = -103.9564799999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 747.5357
$ws.Range("J2").Value = 3112.5
$ws.Range("L2").Value = 3112.5
$ws.Range("N2").Value = -3338.5
$ws.Range("H32").Value = 217812.58
$ws.Range("J32").Value = 15739.6
$ws.Range("L32").Value = 15739.6
$ws.Range("N32").Value = -16313.6
$ws.Range("H97").Value = 7426.75
$ws.Range("J97").Value = 1168.4286
$ws.Range("L97").Value = 1168.4286
$ws.Range("N97").Value = -2160.4286
$ws.Range("H116").Value = 747.5357
$ws.Range("J116").Value = 3112.5
$ws.Range("L116").Value = 3112.5
$ws.Range("N116").Value = -7700.5
$ws.Range("H132").Value = 2991.6287
$ws.Range("I132").Value = 2005.6666
$ws.Range("K132").Value = 6016.9998
$ws.Range("M132").Value = -3486.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 747.5357
$ws.Range("J3").Value = 3112.5
$ws.Range("L3").Value = 3112.5
$ws.Range("N3").Value = -3340.5
$ws.Range("H94").Value = 1747.1333
$ws.Range("I94").Value = 1784.7
$ws.Range("K94").Value = 1784.7
$ws.Range("M94").Value = -1333.7
$ws.Range("H134").Value = 2299.3948
$ws.Range("I134").Value = 1935.5312
$ws.Range("K134").Value = 5806.5936
$ws.Range("M134").Value = -3271.5936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 139.6
$ws.Range("I12").Value = 67.666664
$ws.Range("K12").Value = 67.666664
$ws.Range("M12").Value = 102.333336
$ws.Range("H28").Value = 31927.4
$ws.Range("J28").Value = 31927.4
$ws.Range("L28").Value = 31927.4
$ws.Range("N28").Value = -32417.4
$ws.Range("H31").Value = 3599.2346
$ws.Range("I31").Value = 2290.923
$ws.Range("K31").Value = 2290.923
$ws.Range("M31").Value = -1995.923
$ws.Range("H34").Value = 3599.2346
$ws.Range("I34").Value = 2290.923
$ws.Range("K34").Value = 2290.923
$ws.Range("M34").Value = -2088.923
$ws.Range("H56").Value = 9999.666999999999
$ws.Range("I56").Value = 9999.666999999999
$ws.Range("K56").Value = 9999.666999999999
$ws.Range("M56").Value = -9154.666999999999
$ws.Range("H58").Value = 3446.8215
$ws.Range("J58").Value = 4518.75
$ws.Range("L58").Value = 4518.75
$ws.Range("N58").Value = -4924.75
$ws.Range("H60").Value = 38499.75
$ws.Range("I60").Value = 27499.5
$ws.Range("K60").Value = 27499.5
$ws.Range("M60").Value = -26988.5
$ws.Range("H105").Value = 1934.0476
$ws.Range("I105").Value = 1348.2106
$ws.Range("K105").Value = 1348.2106
$ws.Range("M105").Value = 398.7893999999999
$ws.Range("H107").Value = 2069.4
$ws.Range("I107").Value = 1533.9231
$ws.Range("J107").Value = 5550
$ws.Range("K107").Value = 1533.9231
$ws.Range("L107").Value = 5550
$ws.Range("M107").Value = 386.0769
$ws.Range("N107").Value = -9390
$ws.Range("H132").Value = 12823047
$ws.Range("I132").Value = 2370
$ws.Range("J132").Value = 33336132
$ws.Range("K132").Value = 7110
$ws.Range("L132").Value = 100008396
$ws.Range("M132").Value = -4580
$ws.Range("N132").Value = -100013456
$ws.Range("H136").Value = 3446.8215
$ws.Range("J136").Value = 4518.75
$ws.Range("L136").Value = 13556.25
$ws.Range("N136").Value = -18656.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 681.5
$ws.Range("I2").Value = 454
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 2724
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = -2611
$ws.Range("N2").Value = -6226
$ws.Range("H33").Value = 93361
$ws.Range("I33").Value = 1042.75
$ws.Range("J33").Value = 146114.28
$ws.Range("K33").Value = 6256.5
$ws.Range("L33").Value = 876685.6799999999
$ws.Range("M33").Value = -5973.5
$ws.Range("N33").Value = -877251.6799999999
$ws.Range("H50").Value = 3220.8572
$ws.Range("I50").Value = 2394
$ws.Range("J50").Value = 3358.6667
$ws.Range("K50").Value = 7182
$ws.Range("L50").Value = 10076.0001
$ws.Range("M50").Value = -6701
$ws.Range("N50").Value = -11038.0001
$ws.Range("H53").Value = 3220.8572
$ws.Range("I53").Value = 2394
$ws.Range("J53").Value = 3358.6667
$ws.Range("K53").Value = 7182
$ws.Range("L53").Value = 10076.0001
$ws.Range("M53").Value = -6701
$ws.Range("N53").Value = -11038.0001
$ws.Range("H69").Value = 16337.8
$ws.Range("I69").Value = 50000
$ws.Range("J69").Value = 7922.25
$ws.Range("K69").Value = 150000
$ws.Range("L69").Value = 23766.75
$ws.Range("M69").Value = -149189
$ws.Range("N69").Value = -25388.75
$ws.Range("H72").Value = 16337.8
$ws.Range("I72").Value = 50000
$ws.Range("J72").Value = 7922.25
$ws.Range("K72").Value = 450000
$ws.Range("L72").Value = 71300.25
$ws.Range("M72").Value = -445944
$ws.Range("N72").Value = -79412.25
$ws.Range("H113").Value = 1046.6154
$ws.Range("J113").Value = 1190.875
$ws.Range("L113").Value = 3572.625
$ws.Range("N113").Value = -7912.625
$ws.Range("H122").Value = 3352938.8
$ws.Range("J122").Value = 1136114
$ws.Range("L122").Value = 10225026
$ws.Range("N122").Value = -10229926
$ws.Range("H131").Value = 5489937.5
$ws.Range("I131").Value = 6062277.5
$ws.Range("K131").Value = 18186832.5
$ws.Range("M131").Value = -18181792.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""
$ws.Range("H70").Value = 7102.8433
$ws.Range("I70").Value = 7108.6665
$ws.Range("K70").Value = 7108.6665
$ws.Range("M70").Value = -6838.6665
$ws.Range("H73").Value = 7102.8433
$ws.Range("I73").Value = 7108.6665
$ws.Range("K73").Value = 7108.6665
$ws.Range("M73").Value = -6172.6665
$ws.Range("H97").Value = 583.65515
$ws.Range("I97").Value = 668
$ws.Range("K97").Value = 668
$ws.Range("M97").Value = -172
$ws.Range("H126").Value = 2423.25
$ws.Range("I126").Value = 2398
$ws.Range("K126").Value = 7194
$ws.Range("M126").Value = -4724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2248.3076
$ws.Range("I122").Value = 2244.9048
$ws.Range("K122").Value = 6734.714399999999
$ws.Range("M122").Value = -4284.714399999999
$ws.Range("H136").Value = 3295.9333
$ws.Range("J136").Value = 4286.5
$ws.Range("L136").Value = 12859.5
$ws.Range("N136").Value = -17959.5

Write-Output "applied 223 cell changes"